$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newTimestamp = "2022-07-22 20:57:32"

for ($row = 2; $row -le 73; $row++) {
    $cell = $ws.Cells.Item($row, 15)  # Column O is the 15th column
    $cell.Value = $newTimestamp
}
